# option_book_greeks_data.xlsx -- "Hedging data saving process settled with
# 'export_hedging_data' function"
#
# The exported table gained a row (the hedge function now appends a fresh
# snapshot) and dropped the leading numeric index column: column A becomes
# the "Time" column (previously B), and every other column shifts one to
# the left. The final extent is A1:G10 (header + 9 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grab the existing formatting (bold/bordered header style = the style
# currently on B1:H1; the yyyy-mm-dd hh:mm:ss date style = the one currently
# on B2:B7) before anything is overwritten, and stamp it onto the new
# "Time" column A. Using Copy + PasteSpecial(formats) re-uses the workbook's
# existing cellXfs entries instead of synthesizing new ones.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2:A10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# The old "Time" column (B, rows 2:7) carried the date-style too -- now that
# B is "Spot price" (a plain number), strip that back to the default style
# before writing the new values there (and across the rest of the data
# block, which was already unstyled).
$ws.Range("B2:G10").ClearFormats() | Out-Null

# --- Header row ---
$ws.Range("A1").Value = "Time"
$ws.Range("B1").Value = "Spot price"
$ws.Range("C1").Value = "DELTA leg1"
$ws.Range("D1").Value = "DELTA leg2"
$ws.Range("E1").Value = "DELTA strangle"
$ws.Range("F1").Value = "DELTA hedge"
$ws.Range("G1").Value = "DELTA global"

# --- Data rows (Time, Spot price, DELTA leg1, DELTA leg2, DELTA strangle,
# DELTA hedge, DELTA global) ---
$data = @(
    @(44940.796736111108, 1.07,   -0.49988311216195869, 0.49450767852912231, 0, 0,    0),
    @(44940.796747685177, 1.0828, -0.52351885446815438, 0.47087193622292661, 0, 0.05, 0.05),
    @(44940.798750000002, 1.07,   -0.49988311216195869, 0.49450767852912231, 0, 0,    0),
    @(44940.798784722218, 1.0828, -0.52351885446815438, 0.47087193622292661, 0, 0.05, 0.05),
    @(44940.802442129629, 1.07,   -0.49988311216195869, 0.49450767852912231, 0, 0,    0),
    @(44940.802453703713, 1.0828, -0.52351885446815438, 0.47087193622292661, 0, 0.05, 0.05),
    @(44940.817280092589, 1.07,   -0.49988311216195869, 0.49450767852912231, 0, 0,    0),
    @(44940.817326388889, 1.0828, -0.52351885446815438, 0.47087193622292661, 0, 0.05, 0.05),
    @(44941.45821759259,  1.07,   -0.49988311216195869, 0.49450767852912231, 0, 0,    0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# --- Drop the old trailing column H (used to hold "DELTA global" /
# the old index column's overflow) entirely so the sheet shrinks to G. ---
$ws.Range("H1:H7").Clear() | Out-Null

# --- Column widths: A is wide (32.33203125), B:G share a narrower width
# (13.77734375), replacing the old single-column-B custom width. ---
$ws.Columns.Item(1).ColumnWidth = 32.33203125
$ws.Range("B1:G1").ColumnWidth = 13.77734375

# --- Leave the selection where the export script left the cursor. ---
$ws.Range("H8").Select() | Out-Null
